# Applies the cryptos.xlsx price/volume/hour refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (avoids Excel auto-converting
# numeric-looking strings like "310.92", "1.76%" or "20" into real numbers),
# then restores the default "Normal" style so no stray number-format is left
# behind on the cell (matching the original plain/unstyled data cells).
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "310.92"
Set-TextValue "E2" "1.76%"
Set-TextValue "G2" "20"
# Row 3
Set-TextValue "D3" "35.62"
Set-TextValue "E3" "-1.85%"
Set-TextValue "G3" "20"
# Row 4
Set-TextValue "D4" "5.114"
Set-TextValue "E4" "1.22%"
Set-TextValue "G4" "20"
# Row 5
Set-TextValue "D5" "0.08218"
Set-TextValue "E5" "4.59%"
Set-TextValue "G5" "20"
# Row 6
Set-TextValue "D6" "2.077"
Set-TextValue "E6" "-9.54%"
Set-TextValue "G6" "20"
# Row 7
Set-TextValue "D7" "7.945"
Set-TextValue "E7" "-0.64%"
Set-TextValue "G7" "20"
# Row 8
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D8" "4.128"
Set-TextValue "E8" "-0.77%"
Set-TextValue "G8" "20"
# Row 9
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D9" "2.962"
Set-TextValue "E9" "11.35%"
Set-TextValue "G9" "20"
# Row 10
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D10" "0.9283"
Set-TextValue "E10" "-0.12%"
Set-TextValue "G10" "20"
# Row 11
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D11" "0.1084"
Set-TextValue "E11" "11.13%"
Set-TextValue "G11" "20"
# Row 12
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D12" "0.1927"
Set-TextValue "E12" "3.60%"
Set-TextValue "G12" "20"
# Row 13
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D13" "0.09242"
Set-TextValue "E13" "3.50%"
Set-TextValue "G13" "20"
# Row 14
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D14" "0.03636"
Set-TextValue "E14" "-3.73%"
Set-TextValue "G14" "20"
# Row 15
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D15" "0.09911"
Set-TextValue "E15" "0.13%"
Set-TextValue "G15" "20"
# Row 16
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D16" "0.001436"
Set-TextValue "E16" "-0.81%"
Set-TextValue "G16" "20"
# Row 17
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D17" "0.005774"
Set-TextValue "E17" "1.84%"
Set-TextValue "G17" "20"
# Row 18
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D18" "3.476"
Set-TextValue "E18" "0.41%"
Set-TextValue "G18" "20"
# Row 19
Set-TextValue "E19" "0.33%"
Set-TextValue "G19" "20"
# Row 20
Set-TextValue "D20" "0.1310"
Set-TextValue "E20" "-0.96%"
Set-TextValue "G20" "20"
# Row 21
Set-TextValue "D21" "5.094"
Set-TextValue "E21" "-0.30%"
Set-TextValue "G21" "20"
# Row 22
Set-TextValue "D22" "0.2205"
Set-TextValue "E22" "-2.31%"
Set-TextValue "G22" "20"
# Row 23
Set-TextValue "D23" "0.04547"
Set-TextValue "E23" "-0.73%"
Set-TextValue "G23" "20"
# Row 24
Set-TextValue "E24" "-0.79%"
Set-TextValue "G24" "20"
# Row 25
Set-TextValue "D25" "0.004816"
Set-TextValue "E25" "1.08%"
Set-TextValue "G25" "20"
# Row 26
Set-TextValue "E26" "-4.13%"
Set-TextValue "G26" "20"
# Row 27
Set-TextValue "D27" "0.0004448"
Set-TextValue "E27" "-6.13%"
Set-TextValue "G27" "20"
# Row 28
Set-TextValue "G28" "20"
# Row 29
Set-TextValue "G29" "20"
# Row 30
Set-TextValue "G30" "20"
# Row 31
Set-TextValue "G31" "20"
# Row 32
Set-TextValue "G32" "20"
# Row 33
Set-TextValue "G33" "20"
# Row 34
Set-TextValue "G34" "20"
# Row 35
Set-TextValue "G35" "20"
# Row 36
Set-TextValue "G36" "20"
# Row 37
Set-TextValue "G37" "20"
# Row 38
Set-TextValue "G38" "20"
# Row 39
Set-TextValue "D39" "0.01991"
Set-TextValue "E39" "3.51%"
Set-TextValue "G39" "20"
# Row 40
Set-TextValue "D40" "0.04922"
Set-TextValue "E40" "-3.31%"
Set-TextValue "G40" "20"
# Row 41
Set-TextValue "D41" "0.007838"
Set-TextValue "E41" "-0.18%"
Set-TextValue "G41" "20"
# Row 42
Set-TextValue "D42" "0.009882"
Set-TextValue "E42" "26.02%"
Set-TextValue "G42" "20"
# Row 43
Set-TextValue "D43" "0.1386"
Set-TextValue "E43" "0.02%"
Set-TextValue "G43" "20"
# Row 44
Set-TextValue "E44" "-1.49%"
Set-TextValue "G44" "20"
# Row 45
Set-TextValue "D45" "0.01158"
Set-TextValue "E45" "3.00%"
Set-TextValue "G45" "20"
# Row 46
Set-TextValue "D46" "0.00006541"
Set-TextValue "E46" "5.80%"
Set-TextValue "G46" "20"
# Row 47
Set-TextValue "D47" "0.00000000750"
Set-TextValue "E47" "-0.40%"
Set-TextValue "G47" "20"
# Row 48
Set-TextValue "D48" "177.45"
Set-TextValue "E48" "242.85%"
Set-TextValue "G48" "20"
# Row 49
Set-TextValue "D49" "0.001499"
Set-TextValue "E49" "-21.42%"
Set-TextValue "G49" "20"
# Row 50
Set-TextValue "D50" "0.00002100"
Set-TextValue "E50" "-0.40%"
Set-TextValue "G50" "20"
# Row 51
Set-TextValue "D51" "0.0002000"
Set-TextValue "E51" "-0.40%"
Set-TextValue "G51" "20"
